$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4
$ws.Range("D2").Value = "Caught"
$ws.Range("K2").Value = 39
$ws.Range("L2").Value = 14
$ws.Range("M2").Value = "Caught"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "Bowled"
$ws.Range("E3").Value = " Jasprit Bumrah"
$ws.Range("K3").Value = 10
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = "Caught"
$ws.Range("N3").Value = " Anrich Nortje"
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = "LBW"
$ws.Range("E4").Value = " Bhuvneshwar Kumar"
$ws.Range("K4").Value = 20
$ws.Range("L4").Value = 9
$ws.Range("M4").Value = "Caught"
$ws.Range("B5").Value = 55
$ws.Range("D5").Value = "LBW"
$ws.Range("K5").Value = 43
$ws.Range("L5").Value = 21
$ws.Range("M5").Value = "Caught"
$ws.Range("N5").Value = " Kagiso Rabada"
$ws.Range("B6").Value = 66
$ws.Range("C6").Value = 28
$ws.Range("E6").Value = " Kuldeep Yadav"
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 2
$ws.Range("M6").Value = "Caught"
$ws.Range("N6").Value = " Keshav Maharaj"
$ws.Range("B7").Value = 81
$ws.Range("C7").Value = 28
$ws.Range("D7").Value = "Bowled"
$ws.Range("E7").Value = " Bhuvneshwar Kumar"
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = "Caught"
$ws.Range("N7").Value = " Kagiso Rabada"
$ws.Range("B8").Value = 39
$ws.Range("C8").Value = 17
$ws.Range("D8").Value = "Bowled"
$ws.Range("E8").Value = " Hardik Pandya"
$ws.Range("K8").Value = 14
$ws.Range("L8").Value = 8
$ws.Range("M8").Value = "LBW"
$ws.Range("N8").Value = " Tabraiz Shamsi"
$ws.Range("B9").Value = 21
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = "Caught"
$ws.Range("E9").Value = " Mohommad Shami"
$ws.Range("K9").Value = 37
$ws.Range("L9").Value = 19
$ws.Range("M9").Value = "Caught"
$ws.Range("N9").Value = " Anrich Nortje"
$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = "* NOT OUT"
$ws.Range("E10").Value = " "
$ws.Range("K10").Value = 10
$ws.Range("L10").Value = 6
$ws.Range("B11").Value = 6
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = "NOT OUT"
$ws.Range("E11").Value = " "
$ws.Range("K11").Value = 52
$ws.Range("L11").Value = 13
$ws.Range("M11").Value = "NOT OUT"
$ws.Range("N11").Value = " "
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = " "
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1
$ws.Range("N12").Value = " Anrich Nortje"
$ws.Range("A16").Value = 307
$ws.Range("B16").Value = 8
$ws.Range("C16").Value = "20.0"
$ws.Range("D16").Value = 120
$ws.Range("J16").Value = 233
$ws.Range("L16").Value = "16.3"
$ws.Range("M16").Value = 99
$ws.Range("A21").Value = "Kuldeep Yadav"
$ws.Range("B21").Value = "4.0"
$ws.Range("C21").Value = 53
$ws.Range("E21").Value = 13.25
$ws.Range("J21").Value = "Keshav Maharaj"
$ws.Range("K21").Value = "3.0"
$ws.Range("L21").Value = 49
$ws.Range("M21").Value = 3
$ws.Range("N21").Value = 16.33
$ws.Range("A22").Value = "Jasprit Bumrah"
$ws.Range("B22").Value = "4.0"
$ws.Range("C22").Value = 72
$ws.Range("D22").Value = 2
$ws.Range("J22").Value = "Kagiso Rabada"
$ws.Range("K22").Value = "3.0"
$ws.Range("L22").Value = 42
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 14
$ws.Range("A23").Value = "Mohommad Shami"
$ws.Range("B23").Value = "4.0"
$ws.Range("C23").Value = 62
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 15.5
$ws.Range("J23").Value = "Dwaine Pretorius"
$ws.Range("K23").Value = "3.0"
$ws.Range("L23").Value = 43
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 14.33
$ws.Range("A24").Value = "Bhuvneshwar Kumar"
$ws.Range("B24").Value = "4.0"
$ws.Range("C24").Value = 60
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 15
$ws.Range("J24").Value = "Tabraiz Shamsi"
$ws.Range("K24").Value = "4.0"
$ws.Range("L24").Value = 55
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 13.75
$ws.Range("A25").Value = "Hardik Pandya"
$ws.Range("B25").Value = "4.0"
$ws.Range("C25").Value = 60
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 15
$ws.Range("J25").Value = "Anrich Nortje"
$ws.Range("K25").Value = "3.3"
$ws.Range("L25").Value = 44
$ws.Range("M25").Value = 3
$ws.Range("N25").Value = 13.33
